# Scheduled runner update: refresh computed market-price/profit figures
# (currentAveragePrice*, LevePrice*, LeveProfit* columns H-N) on several
# leve rows across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 667.0345
$ws.Range("J17").Value = 676.5714
$ws.Range("L17").Value = 2029.7142
$ws.Range("N17").Value = -2365.7142

$ws.Range("H98").Value = 1020363.44
$ws.Range("I98").Value = 1245999.8
$ws.Range("K98").Value = 1245999.8
$ws.Range("M98").Value = -1244501.8

$ws.Range("H116").Value = 7134716.5
$ws.Range("I116").Value = 11669589
$ws.Range("J116").Value = 8487.429
$ws.Range("K116").Value = 11669589
$ws.Range("L116").Value = 8487.429
$ws.Range("M116").Value = -11666147
$ws.Range("N116").Value = -15371.429

$ws.Range("H122").Value = 1020363.44
$ws.Range("I122").Value = 1245999.8
$ws.Range("K122").Value = 3737999.4
$ws.Range("M122").Value = -3735549.4

$ws.Range("H125").Value = 12457301
$ws.Range("I125").Value = 599
$ws.Range("J125").Value = 14014388
$ws.Range("K125").Value = 5391
$ws.Range("L125").Value = 126129492
$ws.Range("M125").Value = -2931
$ws.Range("N125").Value = -126134412

$ws.Range("H137").Value = 23257016
$ws.Range("I137").Value = 32258916
$ws.Range("J137").Value = 2109.6667
$ws.Range("K137").Value = 96776748
$ws.Range("L137").Value = 6329.000100000001
$ws.Range("M137").Value = -96774198
$ws.Range("N137").Value = -11429.0001


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1405.5834
$ws.Range("I2").Value = 932.8333
$ws.Range("J2").Value = 1878.3334
$ws.Range("K2").Value = 932.8333
$ws.Range("L2").Value = 1878.3334
$ws.Range("M2").Value = -819.8333
$ws.Range("N2").Value = -2104.3334

$ws.Range("H32").Value = 4293.755
$ws.Range("I32").Value = 1992.0605
$ws.Range("K32").Value = 1992.0605
$ws.Range("M32").Value = -1705.0605

$ws.Range("H116").Value = 1405.5834
$ws.Range("I116").Value = 932.8333
$ws.Range("J116").Value = 1878.3334
$ws.Range("K116").Value = 932.8333
$ws.Range("L116").Value = 1878.3334
$ws.Range("M116").Value = 1361.1667
$ws.Range("N116").Value = -6466.3334

$ws.Range("H122").Value = 2711.0833
$ws.Range("I122").Value = 2837.2222
$ws.Range("K122").Value = 8511.6666
$ws.Range("M122").Value = -6061.6666

$ws.Range("H132").Value = 2333.7812
$ws.Range("I132").Value = 2030.1923
$ws.Range("K132").Value = 6090.5769
$ws.Range("M132").Value = -3560.5769


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1405.5834
$ws.Range("I3").Value = 932.8333
$ws.Range("J3").Value = 1878.3334
$ws.Range("K3").Value = 932.8333
$ws.Range("L3").Value = 1878.3334
$ws.Range("M3").Value = -818.8333
$ws.Range("N3").Value = -2106.3334

$ws.Range("H134").Value = 5316.0713
$ws.Range("I134").Value = 4655.875
$ws.Range("K134").Value = 13967.625
$ws.Range("M134").Value = -11432.625


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1240.4642
$ws.Range("I31").Value = 1137.32
$ws.Range("J31").Value = 2100
$ws.Range("K31").Value = 1137.32
$ws.Range("L31").Value = 2100
$ws.Range("M31").Value = -842.3199999999999
$ws.Range("N31").Value = -2690

$ws.Range("H34").Value = 1240.4642
$ws.Range("I34").Value = 1137.32
$ws.Range("J34").Value = 2100
$ws.Range("K34").Value = 1137.32
$ws.Range("L34").Value = 2100
$ws.Range("M34").Value = -935.3199999999999
$ws.Range("N34").Value = -2504

$ws.Range("H56").Value = 8950.5
$ws.Range("I56").Value = 6999.6665
$ws.Range("J56").Value = 14803
$ws.Range("K56").Value = 6999.6665
$ws.Range("L56").Value = 14803
$ws.Range("M56").Value = -6154.6665
$ws.Range("N56").Value = -16493

$ws.Range("H86").Value = 62501910
$ws.Range("J86").Value = 2130
$ws.Range("L86").Value = 2130
$ws.Range("N86").Value = -4376

$ws.Range("H89").Value = 62501910
$ws.Range("J89").Value = 2130
$ws.Range("L89").Value = 10650
$ws.Range("N89").Value = -21882

$ws.Range("H122").Value = 2148.6316
$ws.Range("I122").Value = 1489.8889
$ws.Range("J122").Value = 2741.5
$ws.Range("K122").Value = 4469.6667
$ws.Range("L122").Value = 8224.5
$ws.Range("M122").Value = -2019.6667
$ws.Range("N122").Value = -13124.5


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1298.4584
$ws.Range("I5").Value = 980.8
$ws.Range("K5").Value = 2942.4
$ws.Range("M5").Value = -2830.4

$ws.Range("H113").Value = 19231670
$ws.Range("J113").Value = 20000912
$ws.Range("L113").Value = 60002736
$ws.Range("N113").Value = -60007076

$ws.Range("H135").Value = 1298.4584
$ws.Range("I135").Value = 980.8
$ws.Range("K135").Value = 8827.199999999999
$ws.Range("M135").Value = -6292.199999999999


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2406.2856
$ws.Range("I80").Value = 2072
$ws.Range("J80").Value = 2540
$ws.Range("K80").Value = 2072
$ws.Range("L80").Value = 2540
$ws.Range("M80").Value = -1074
$ws.Range("N80").Value = -4536

$ws.Range("H83").Value = 2406.2856
$ws.Range("I83").Value = 2072
$ws.Range("J83").Value = 2540
$ws.Range("K83").Value = 10360
$ws.Range("L83").Value = 12700
$ws.Range("M83").Value = -5368
$ws.Range("N83").Value = -22684

$ws.Range("H122").Value = 2223244.2
$ws.Range("I122").Value = 3704236.8
$ws.Range("J122").Value = 1755.5
$ws.Range("K122").Value = 11112710.4
$ws.Range("L122").Value = 5266.5
$ws.Range("M122").Value = -11110260.4
$ws.Range("N122").Value = -10166.5

$ws.Range("H126").Value = 2570.9583
$ws.Range("I126").Value = 1732.4445
$ws.Range("J126").Value = 3074.0667
$ws.Range("K126").Value = 5197.333500000001
$ws.Range("L126").Value = 9222.2001
$ws.Range("M126").Value = -2727.333500000001
$ws.Range("N126").Value = -14162.2001

$ws.Range("H132").Value = 2783
$ws.Range("I132").Value = 2643.1738
$ws.Range("J132").Value = 5999
$ws.Range("K132").Value = 7929.5214
$ws.Range("L132").Value = 17997
$ws.Range("M132").Value = -5399.5214
$ws.Range("N132").Value = -23057


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3357.1428
$ws.Range("I40").Value = 1500
$ws.Range("J40").Value = 3500
$ws.Range("K40").Value = 1500
$ws.Range("L40").Value = 3500
$ws.Range("M40").Value = -1364
$ws.Range("N40").Value = -3772

$ws.Range("H46").Value = 528.5
$ws.Range("I46").Value = 508.42856
$ws.Range("J46").Value = 548.5714
$ws.Range("K46").Value = 508.42856
$ws.Range("L46").Value = 548.5714
$ws.Range("M46").Value = -320.42856
$ws.Range("N46").Value = -924.5714


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 550
$ws.Range("I113").Value = 300
$ws.Range("J113").Value = 800
$ws.Range("K113").Value = 900
$ws.Range("L113").Value = 2400
$ws.Range("M113").Value = 1270
$ws.Range("N113").Value = -6740

